$d = $word.ActiveDocument

# Helper: returns a Range covering a paragraph's text but excluding its
# trailing paragraph mark, so that Font formatting only lands on the
# run(s) inside <w:r>, not on the paragraph mark's own run properties
# (<w:pPr><w:rPr>).
function Get-ParaTextRange($para) {
    $r = $para.Range
    return $d.Range($r.Start, $r.End - 1)
}

# 1) Title paragraph ("Uplift Engine 2.1 — Project Summary" + line break,
#    centered, bold): add Times New Roman font and bump size 24pt -> 26pt
#    (w:sz 48 -> 52, stored in half-points).
$pTitle = $d.Paragraphs.Item(1)
if ($pTitle.Range.Text -like "Uplift Engine*Project Summary*") {
    $rTitleRun = Get-ParaTextRange $pTitle
    $rTitleRun.Font.Name = "Times New Roman"
    $rTitleRun.Font.Size = 26
}

# 2) Author/Date paragraph on the title page (centered, "Author: Team
#    Uplift" / "Date: ..."): add Arial 11pt (w:sz 22).
$pAuthor = $d.Paragraphs.Item(3)
if ($pAuthor.Range.Text -like "Author: Team Uplift*Date:*") {
    $rAuthorRun = Get-ParaTextRange $pAuthor
    $rAuthorRun.Font.Name = "Arial"
    $rAuthorRun.Font.Size = 11
}

# 3) Heading1 paragraph ("Uplift Engine 2.1 — Project Summary"): add
#    Times New Roman 18pt (w:sz 36).
$pHeading = $d.Paragraphs.Item(7)
if ($pHeading.Range.Text -like "Uplift Engine*Project Summary*") {
    $rHeadingRun = Get-ParaTextRange $pHeading
    $rHeadingRun.Font.Name = "Times New Roman"
    $rHeadingRun.Font.Size = 18
}

# 4) Normal paragraph style default run formatting: Arial 11pt (w:sz 22).
#    This is the md->docx generator's body-text default.
$normalStyle = $d.Styles("Normal")
$normalStyle.Font.Name = "Arial"
$normalStyle.Font.Size = 11
